$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d = $ws.Range("D2")
$d.NumberFormat = "@"
$d.Value = "27.660.63"
$d.Style = "Normal"
$e = $ws.Range("E2")
$e.NumberFormat = "@"
$e.Value = "  +1.06%  "
$e.Style = "Normal"
$d = $ws.Range("D3")
$d.NumberFormat = "@"
$d.Value = "1.871.10"
$d.Style = "Normal"
$e = $ws.Range("E3")
$e.NumberFormat = "@"
$e.Value = "  +0.47%  "
$e.Style = "Normal"
$d = $ws.Range("D4")
$d.NumberFormat = "@"
$d.Value = "1.003"
$d.Style = "Normal"
$e = $ws.Range("E4")
$e.NumberFormat = "@"
$e.Value = "  +0.34%  "
$e.Style = "Normal"
$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = "332.12"
$d.Style = "Normal"
$e = $ws.Range("E5")
$e.NumberFormat = "@"
$e.Value = "  +2.60%  "
$e.Style = "Normal"
$e = $ws.Range("E6")
$e.NumberFormat = "@"
$e.Value = "  +0.34%  "
$e.Style = "Normal"
$d = $ws.Range("D7")
$d.NumberFormat = "@"
$d.Value = "0.4731"
$d.Style = "Normal"
$e = $ws.Range("E7")
$e.NumberFormat = "@"
$e.Value = "  +4.36%  "
$e.Style = "Normal"
$d = $ws.Range("D8")
$d.NumberFormat = "@"
$d.Value = "0.3941"
$d.Style = "Normal"
$e = $ws.Range("E8")
$e.NumberFormat = "@"
$e.Value = "  +1.74%  "
$e.Style = "Normal"
$d = $ws.Range("D9")
$d.NumberFormat = "@"
$d.Value = "47.92"
$d.Style = "Normal"
$e = $ws.Range("E9")
$e.NumberFormat = "@"
$e.Value = "  -0.54%  "
$e.Style = "Normal"
$d = $ws.Range("D10")
$d.NumberFormat = "@"
$d.Value = "0.08054"
$d.Style = "Normal"
$e = $ws.Range("E10")
$e.NumberFormat = "@"
$e.Value = "  +1.73%  "
$e.Style = "Normal"
$d = $ws.Range("D11")
$d.NumberFormat = "@"
$d.Value = "1.025"
$d.Style = "Normal"
$e = $ws.Range("E11")
$e.NumberFormat = "@"
$e.Value = "  +0.25%  "
$e.Style = "Normal"
$d = $ws.Range("D12")
$d.NumberFormat = "@"
$d.Value = "21.99"
$d.Style = "Normal"
$e = $ws.Range("E12")
$e.NumberFormat = "@"
$e.Value = "  +2.61%  "
$e.Style = "Normal"
$d = $ws.Range("D13")
$d.NumberFormat = "@"
$d.Value = "1.890.53"
$d.Style = "Normal"
$e = $ws.Range("E13")
$e.NumberFormat = "@"
$e.Value = "  +1.55%  "
$e.Style = "Normal"
$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = "5.950"
$d.Style = "Normal"
$e = $ws.Range("E14")
$e.NumberFormat = "@"
$e.Value = "  +0.72%  "
$e.Style = "Normal"
$d = $ws.Range("D15")
$d.NumberFormat = "@"
$d.Value = "7.141"
$d.Style = "Normal"
$e = $ws.Range("E15")
$e.NumberFormat = "@"
$e.Value = "  -0.19%  "
$e.Style = "Normal"
$d = $ws.Range("D16")
$d.NumberFormat = "@"
$d.Value = "1.006"
$d.Style = "Normal"
$e = $ws.Range("E16")
$e.NumberFormat = "@"
$e.Value = "  +0.60%  "
$e.Style = "Normal"
$e = $ws.Range("E17")
$e.NumberFormat = "@"
$e.Value = "  +1.24%  "
$e.Style = "Normal"
$d = $ws.Range("D18")
$d.NumberFormat = "@"
$d.Value = "86.82"
$d.Style = "Normal"
$e = $ws.Range("E18")
$e.NumberFormat = "@"
$e.Value = "  +1.18%  "
$e.Style = "Normal"
$d = $ws.Range("D19")
$d.NumberFormat = "@"
$d.Value = "0.06653"
$d.Style = "Normal"
$e = $ws.Range("E19")
$e.NumberFormat = "@"
$e.Value = "  +2.28%  "
$e.Style = "Normal"
$d = $ws.Range("D20")
$d.NumberFormat = "@"
$d.Value = "17.12"
$d.Style = "Normal"
$e = $ws.Range("E20")
$e.NumberFormat = "@"
$e.Value = "  -0.05%  "
$e.Style = "Normal"
$e = $ws.Range("E21")
$e.NumberFormat = "@"
$e.Value = "  +0.29%  "
$e.Style = "Normal"
$d = $ws.Range("D22")
$d.NumberFormat = "@"
$d.Value = "27.665.78"
$d.Style = "Normal"
$e = $ws.Range("E22")
$e.NumberFormat = "@"
$e.Value = "  +1.13%  "
$e.Style = "Normal"
$d = $ws.Range("D23")
$d.NumberFormat = "@"
$d.Value = "5.500"
$d.Style = "Normal"
$e = $ws.Range("E23")
$e.NumberFormat = "@"
$e.Value = "  -0.66%  "
$e.Style = "Normal"
$d = $ws.Range("D24")
$d.NumberFormat = "@"
$d.Value = "10.98"
$d.Style = "Normal"
$e = $ws.Range("E24")
$e.NumberFormat = "@"
$e.Value = "  +0.78%  "
$e.Style = "Normal"
$d = $ws.Range("D25")
$d.NumberFormat = "@"
$d.Value = "2.307"
$d.Style = "Normal"
$e = $ws.Range("E25")
$e.NumberFormat = "@"
$e.Value = "  +1.32%  "
$e.Style = "Normal"
$d = $ws.Range("D26")
$d.NumberFormat = "@"
$d.Value = "2.098.44"
$d.Style = "Normal"
$e = $ws.Range("E26")
$e.NumberFormat = "@"
$e.Value = "  +0.88%  "
$e.Style = "Normal"
$d = $ws.Range("D27")
$d.NumberFormat = "@"
$d.Value = "158.84"
$d.Style = "Normal"
$e = $ws.Range("E27")
$e.NumberFormat = "@"
$e.Value = "  +3.88%  "
$e.Style = "Normal"
$d = $ws.Range("D28")
$d.NumberFormat = "@"
$d.Value = "20.28"
$d.Style = "Normal"
$e = $ws.Range("E28")
$e.NumberFormat = "@"
$e.Value = "  +2.33%  "
$e.Style = "Normal"
$d = $ws.Range("D29")
$d.NumberFormat = "@"
$d.Value = "2.097"
$d.Style = "Normal"
$e = $ws.Range("E29")
$e.NumberFormat = "@"
$e.Value = "  +1.55%  "
$e.Style = "Normal"
$d = $ws.Range("D30")
$d.NumberFormat = "@"
$d.Value = "5.556"
$d.Style = "Normal"
$e = $ws.Range("E30")
$e.NumberFormat = "@"
$e.Value = "  +1.13%  "
$e.Style = "Normal"
$e = $ws.Range("E31")
$e.NumberFormat = "@"
$e.Value = "  +1.04%  "
$e.Style = "Normal"
$d = $ws.Range("D32")
$d.NumberFormat = "@"
$d.Value = "0.9700"
$d.Style = "Normal"
$e = $ws.Range("E32")
$e.NumberFormat = "@"
$e.Value = "  +3.70%  "
$e.Style = "Normal"
$d = $ws.Range("D33")
$d.NumberFormat = "@"
$d.Value = "0.09530"
$d.Style = "Normal"
$e = $ws.Range("E33")
$e.NumberFormat = "@"
$e.Value = "  +2.15%  "
$e.Style = "Normal"
$d = $ws.Range("D34")
$d.NumberFormat = "@"
$d.Value = "1.450"
$d.Style = "Normal"
$e = $ws.Range("E34")
$e.NumberFormat = "@"
$e.Value = "  -3.24%  "
$e.Style = "Normal"
$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = "3.588"
$d.Style = "Normal"
$e = $ws.Range("E35")
$e.NumberFormat = "@"
$e.Value = "  -0.52%  "
$e.Style = "Normal"
$d = $ws.Range("D36")
$d.NumberFormat = "@"
$d.Value = "5.332"
$d.Style = "Normal"
$e = $ws.Range("E36")
$e.NumberFormat = "@"
$e.Value = "  +1.01%  "
$e.Style = "Normal"
$d = $ws.Range("D37")
$d.NumberFormat = "@"
$d.Value = "0.06103"
$d.Style = "Normal"
$e = $ws.Range("E37")
$e.NumberFormat = "@"
$e.Value = "  +1.65%  "
$e.Style = "Normal"
$d = $ws.Range("D38")
$d.NumberFormat = "@"
$d.Value = "0.02260"
$d.Style = "Normal"
$e = $ws.Range("E38")
$e.NumberFormat = "@"
$e.Value = "  +0.93%  "
$e.Style = "Normal"
$d = $ws.Range("D39")
$d.NumberFormat = "@"
$d.Value = "1.224"
$d.Style = "Normal"
$e = $ws.Range("E39")
$e.NumberFormat = "@"
$e.Value = "  +0.08%  "
$e.Style = "Normal"
$d = $ws.Range("D40")
$d.NumberFormat = "@"
$d.Value = "8.165"
$d.Style = "Normal"
$e = $ws.Range("E40")
$e.NumberFormat = "@"
$e.Value = "  -1.09%  "
$e.Style = "Normal"
$d = $ws.Range("D41")
$d.NumberFormat = "@"
$d.Value = "0.6032"
$d.Style = "Normal"
$e = $ws.Range("E41")
$e.NumberFormat = "@"
$e.Value = "  +1.97%  "
$e.Style = "Normal"
$e = $ws.Range("E42")
$e.NumberFormat = "@"
$e.Value = "  +0.34%  "
$e.Style = "Normal"
$d = $ws.Range("D43")
$d.NumberFormat = "@"
$d.Value = "10.26"
$d.Style = "Normal"
$e = $ws.Range("E43")
$e.NumberFormat = "@"
$e.Value = "  +1.09%  "
$e.Style = "Normal"
$e = $ws.Range("E44")
$e.NumberFormat = "@"
$e.Value = "  -1.37%  "
$e.Style = "Normal"
$d = $ws.Range("D45")
$d.NumberFormat = "@"
$d.Value = "0.5715"
$d.Style = "Normal"
$e = $ws.Range("E45")
$e.NumberFormat = "@"
$e.Value = "  +1.46%  "
$e.Style = "Normal"
$d = $ws.Range("D46")
$d.NumberFormat = "@"
$d.Value = "12.23"
$d.Style = "Normal"
$e = $ws.Range("E46")
$e.NumberFormat = "@"
$e.Value = "  +1.81%  "
$e.Style = "Normal"
$d = $ws.Range("D47")
$d.NumberFormat = "@"
$d.Value = "1.942"
$d.Style = "Normal"
$e = $ws.Range("E47")
$e.NumberFormat = "@"
$e.Value = "  +0.65%  "
$e.Style = "Normal"
$d = $ws.Range("D48")
$d.NumberFormat = "@"
$d.Value = "3.386"
$d.Style = "Normal"
$e = $ws.Range("E48")
$e.NumberFormat = "@"
$e.Value = "  +0.45%  "
$e.Style = "Normal"
$d = $ws.Range("D49")
$d.NumberFormat = "@"
$d.Value = "0.06863"
$d.Style = "Normal"
$e = $ws.Range("E49")
$e.NumberFormat = "@"
$e.Value = "  +1.06%  "
$e.Style = "Normal"
$d = $ws.Range("D50")
$d.NumberFormat = "@"
$d.Value = "114.57"
$d.Style = "Normal"
$e = $ws.Range("E50")
$e.NumberFormat = "@"
$e.Value = "  +5.96%  "
$e.Style = "Normal"
$d = $ws.Range("D51")
$d.NumberFormat = "@"
$d.Value = "0.00000000302"
$d.Style = "Normal"
$e = $ws.Range("E51")
$e.NumberFormat = "@"
$e.Value = "  +15.97%  "
$e.Style = "Normal"
